$d = $word.ActiveDocument

# The "Notice u/s 94 BNSS, 2023" heading belongs right before the "To," salutation
# paragraph (after the three blank letterhead lines). Locate "To," by text so the
# insertion point is correct even if the template layout shifts.
$toPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "To,`r") {
        $toPara = $d.Paragraphs($i)
        break
    }
}

# Create a brand-new (initially empty) paragraph immediately before "To,". After this
# call, $toPara's Range now refers to that freshly inserted empty paragraph (the "To,"
# text itself shifted down into the next paragraph).
$toPara.Range.InsertParagraphBefore()

# Overwrite the new empty paragraph's contents with the exact OOXML for the
# centered, bold, single-underlined notice heading. Using raw XML (rather than
# Range.Text + Font.* property assignments) avoids inheriting the surrounding
# "Body A" paragraph/run formatting (fonts, sizes, spacing, justification, etc.).
$fragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Notice u/s 94 BNSS, 2023</w:t></w:r></w:p>'
$toPara.Range.InsertXML($fragment)
